$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "38.220.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.53%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.059.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.25%  "

# Row 4
$ws.Range("E4").Value = "  +0.20%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.02%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.616"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.02%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.60"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.10%  "

# Row 8
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
$ws.Range("E9").Value = "  +3.70%  "

# Row 10
$ws.Range("E10").Value = "  +4.18%  "

# Row 11
$ws.Range("E11").Value = "  +0.73%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.364.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.30%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.69"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.87%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.78%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.756"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.10%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.74%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.057.95"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.38%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "38.090.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.37%  "

# Row 19
$ws.Range("E19").Value = "  +1.54%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.37%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0833"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.24%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "225.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.50%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.03%  "

# Row 24
$ws.Range("E24").Value = "  +0.92%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.35%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.97%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.77%  "

# Row 28
$ws.Range("E28").Value = "  +8.96%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.80%  "

# Row 30
$ws.Range("E30").Value = "  +2.28%  "

# Row 31
$ws.Range("E31").Value = "  +2.65%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.56"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.20%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.63"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.10%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0615"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.78%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.98"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.44%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.34"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.39%  "

# Row 37
$ws.Range("E37").Value = "  +15.86%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.86%  "

# Row 39
$ws.Range("E39").Value = "  -0.16%  "

# Row 40
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0220"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.06%  "

# Row 41
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.81%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.481.83"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.54%  "

# Row 43
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0948"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.93%  "

# Row 44
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.89"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.34%  "

# Row 45
$ws.Range("E45").Value = "  +4.29%  "

# Row 46
$ws.Range("E46").Value = "  +1.30%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +17.51%  "

# Row 48
$ws.Range("E48").Value = "  +2.42%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.47%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.10"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.16%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.250.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.27%  "
